{"js": "const replacements = [\n  [\"2023-11-04 Saturday\", \"2023-11-05 Sunday\"],\n  [\"67\u00d768=4556\", \"28\u00d791=2548\"],\n  [\"81\u00d769=5589\", \"62\u00d739=2418\"],\n  [\"68\u00d738=2584\", \"18\u00d731=558\"],\n  [\"48\u00d741=1968\", \"99\u00d720=1980\"],\n  [\"93\u00d792=8556\", \"34\u00d754=1836\"],\n  [\"72\u00d736=2592\", \"60\u00d752=3120\"],\n  [\"15\u00d790=1350\", \"89\u00d728=2492\"],\n  [\"54\u00d792=4968\", \"42\u00d751=2142\"],\n  [\"75\u00d772=5400\", \"64\u00d741=2624\"],\n  [\"32\u00d713=416\", \"87\u00d738=3306\"],\n  [\"54\u00d777=4158\", \"95\u00d713=1235\"],\n  [\"86\u00d745=3870\", \"95\u00d751=4845\"],\n  [\"25\u00d755=1375\", \"71\u00d719=1349\"],\n  [\"40\u00d717=680\", \"87\u00d739=3393\"],\n  [\"24\u00d769=1656\", \"34\u00d782=2788\"],\n  [\"68\u00d756=3808\", \"59\u00d739=2301\"],\n  [\"75\u00d750=3750\", \"72\u00d775=5400\"],\n  [\"51\u00d757=2907\", \"31\u00d720=620\"],\n  [\"74\u00d741=3034\", \"28\u00d721=588\"],\n  [\"77\u00d751=3927\", \"50\u00d731=1550\"],\n  [\"57\u00d738=2166\", \"62\u00d774=4588\"],\n  [\"46\u00d763=2898\", \"26\u00d753=1378\"],\n  [\"61\u00d783=5063\", \"34\u00d712=408\"],\n  [\"30\u00d789=2670\", \"52\u00d789=4628\"],\n  [\"18\u00d724=432\", \"23\u00d777=1771\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$pairs = @(\n    @(\"2023-11-04 Saturday\", \"2023-11-05 Sunday\"),\n    @(\"67\u00d768=4556\", \"28\u00d791=2548\"),\n    @(\"81\u00d769=5589\", \"62\u00d739=2418\"),\n    @(\"68\u00d738=2584\", \"18\u00d731=558\"),\n    @(\"48\u00d741=1968\", \"99\u00d720=1980\"),\n    @(\"93\u00d792=8556\", \"34\u00d754=1836\"),\n    @(\"72\u00d736=2592\", \"60\u00d752=3120\"),\n    @(\"15\u00d790=1350\", \"89\u00d728=2492\"),\n    @(\"54\u00d792=4968\", \"42\u00d751=2142\"),\n    @(\"75\u00d772=5400\", \"64\u00d741=2624\"),\n    @(\"32\u00d713=416\", \"87\u00d738=3306\"),\n    @(\"54\u00d777=4158\", \"95\u00d713=1235\"),\n    @(\"86\u00d745=3870\", \"95\u00d751=4845\"),\n    @(\"25\u00d755=1375\", \"71\u00d719=1349\"),\n    @(\"40\u00d717=680\", \"87\u00d739=3393\"),\n    @(\"24\u00d769=1656\", \"34\u00d782=2788\"),\n    @(\"68\u00d756=3808\", \"59\u00d739=2301\"),\n    @(\"75\u00d750=3750\", \"72\u00d775=5400\"),\n    @(\"51\u00d757=2907\", \"31\u00d720=620\"),\n    @(\"74\u00d741=3034\", \"28\u00d721=588\"),\n    @(\"77\u00d751=3927\", \"50\u00d731=1550\"),\n    @(\"57\u00d738=2166\", \"62\u00d774=4588\"),\n    @(\"46\u00d763=2898\", \"26\u00d753=1378\"),\n    @(\"61\u00d783=5063\", \"34\u00d712=408\"),\n    @(\"30\u00d789=2670\", \"52\u00d789=4628\"),\n    @(\"18\u00d724=432\", \"23\u00d777=1771\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}"}
